$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AF4").Value = 0.783
$ws.Range("AF5").Value = 0.95
$ws.Range("AF6").Value = 0.858
$ws.Range("AF7").Value = 0.911
$ws.Range("AF8").Value = 0.888
$ws.Range("AF9").Value = 0.8
$ws.Range("AF10").Value = 0.95
$ws.Range("AF11").Value = 0.95
$ws.Range("AF12").Value = 1.211
$ws.Range("AF13").Value = 1.4
